$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.042.25"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.645.16"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "1.874.98"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.636.02"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.531"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "27.019.83"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.118"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0513"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("D35").Value = "1.265.46"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0173"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.534"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "1.785.55"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  -7.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
